# Commit: "Model fitting based on data preprocessing"
# The traffic_volume (column D) values on the "Training Data" sheet are
# replaced with re-derived/pre-processed figures, and the (previously
# embedded, very large) raw per-record inline-string data dumps in
# column D of the "Testing Data" sheet are cleared out.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Training Data sheet: update traffic_volume (column D) for rows 2-61
# ---------------------------------------------------------------------
$wsTrain = $wb.Worksheets.Item("Training Data")

$newTrafficVolume = @(
    14008,10402,10825,10366,9232,9406,12776,10709,10444,10484,
    10366,15015,14115,9357,10957,10556,10108,10134,13486,11162,
    10535,11038,10661,15325,14947,10907,9399,7543,8943,9327,
    10249,10424,10894,11766,11016,15255,13657,351,278,309,
    327,331,384,357,366,370,357,471,410,336,
    333,359,323,328,403,359,365,364,357,522
)

$startRow = 2
for ($i = 0; $i -lt $newTrafficVolume.Length; $i++) {
    $row = $startRow + $i
    $wsTrain.Cells.Item($row, 4).Value = $newTrafficVolume[$i]
}

# ---------------------------------------------------------------------
# 2) Testing Data sheet: clear the bulky inline-string dumps that were
#    stored in column D for rows 2-13
# ---------------------------------------------------------------------
$wsTest = $wb.Worksheets.Item("Testing Data")
$wsTest.Range("D2:D13").ClearContents()
